# Updated 2D training schedules, no break screen
# Adds a new "break_on_off" column (L) to Sheet1, with a header in L1 and
# a 0/1 flag per trial row (rows 2-73).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column header (becomes shared-string index 21: "break_on_off")
$ws.Range("L1").Value = "break_on_off"

# Per-trial break_on_off flags for rows 2..73 (trial 1..72)
$breakOnOff = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $breakOnOff.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $breakOnOff[$i]
}

# Update the view: select the newly added column instead of the old M35 cell
[void]$ws.Range("L1:L73").Select()
